$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D) text, new Volume(1h) (column E) text.
# "" (empty string) means that column is unchanged for that row.
$updates = @(
    ,@(2, "27.374.99", "-2.84%")
    ,@(3, "1.739.98", "-3.52%")
    ,@(4, "1.003", "+0.06%")
    ,@(5, "322.91", "-3.91%")
    ,@(6, "1.000", "+0.04%")
    ,@(7, "0.4238", "-9.27%")
    ,@(8, "0.3612", "-2.78%")
    ,@(9, "45.47", "+0.20%")
    ,@(10, "0.07408", "-3.64%")
    ,@(11, "1.112", "-3.66%")
    ,@(12, "1.002", "+0.04%")
    ,@(13, "21.51", "-5.00%")
    ,@(14, "", "-4.96%")
    ,@(15, "7.173", "-3.17%")
    ,@(16, "1.730.95", "-3.77%")
    ,@(17, "0.00001061", "-3.19%")
    ,@(18, "87.40", "+5.76%")
    ,@(19, "0.06004", "-10.86%")
    ,@(20, "", "-0.02%")
    ,@(21, "16.81", "-3.80%")
    ,@(22, "6.079", "-5.55%")
    ,@(23, "0.5227", "-5.96%")
    ,@(24, "27.375.17", "-2.77%")
    ,@(25, "", "-4.97%")
    ,@(26, "2.382", "-1.20%")
    ,@(27, "20.12", "-3.78%")
    ,@(28, "2.366", "-1.44%")
    ,@(29, "149.14", "-2.41%")
    ,@(30, "1.927.75", "-3.91%")
    ,@(31, "126.37", "-6.13%")
    ,@(32, "1.187", "-6.57%")
    ,@(33, "5.660", "-4.51%")
    ,@(34, "0.09093", "-5.87%")
    ,@(35, "3.645", "-9.86%")
    ,@(36, "12.91", "+5.42%")
    ,@(37, "0.2134", "-5.63%")
    ,@(38, "5.062", "-3.96%")
    ,@(39, "0.02249", "-5.38%")
    ,@(40, "0.06049", "-5.62%")
    ,@(41, "0.6371", "-5.48%")
    ,@(42, "1.185", "-4.16%")
    ,@(43, "7.958", "-2.14%")
    ,@(44, "", "+0.00%")
    ,@(45, "1.407", "-7.57%")
    ,@(46, "13.57", "")
    ,@(47, "3.719", "-3.18%")
    ,@(48, "0.5818", "-6.17%")
    ,@(49, "125.02", "-3.85%")
    ,@(50, "1.955", "-5.62%")
    ,@(51, "0.06849", "-4.18%")
)

foreach ($u in $updates) {
    $row  = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]

    if ($dVal -ne "") {
        # Force text storage so values like "1.003" or "27.374.99" are not
        # reinterpreted by Excel as numbers, then restore the default style
        # so no stray number-format styling is left on the cell.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    }

    if ($eVal -ne "") {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = "  " + $eVal + "  "
        $eCell.Style = "Normal"
    }
}
